$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly measured headspace preparation rows (Summer 2021 data continuation)
$ws.Range("A14").Value = 44341
$ws.Range("B14").Value = 20.9
$ws.Range("C14").Value = 30.18

$ws.Range("A15").Value = 44348
$ws.Range("B15").Value = 20.7
$ws.Range("C15").Value = 30.242999999999999

$ws.Range("A16").Value = 44361
$ws.Range("B16").Value = 20.7
$ws.Range("C16").Value = 29.858000000000001

$ws.Range("A17").Value = 44364
$ws.Range("B17").Value = 20.6
$ws.Range("C17").Value = 30.053999999999998

$ws.Range("A18").Value = 44370
$ws.Range("B18").Value = 21
$ws.Range("C18").Value = 30.184000000000001

$ws.Range("A19").Value = 44376
$ws.Range("B19").Value = 23.2
$ws.Range("C19").Value = 30.222000000000001

$ws.Range("A20").Value = 44386
$ws.Range("B20").Value = 22.5
$ws.Range("C20").Value = 29.928999999999998

$ws.Range("A21").Value = 44390
$ws.Range("B21").Value = 23.8
$ws.Range("C21").Value = 30.183

$ws.Range("A22").Value = 44399
$ws.Range("B22").Value = 22.3
$ws.Range("C22").Value = 30.111000000000001

$ws.Range("A23").Value = 44404
$ws.Range("B23").Value = 21.9
$ws.Range("C23").Value = 29.977

$ws.Range("A24").Value = 44411
$ws.Range("B24").Value = 22.1
$ws.Range("C24").Value = 30.03

$ws.Range("A25").Value = 44418
$ws.Range("B25").Value = 22.2
$ws.Range("C25").Value = 30.135999999999999

$ws.Range("A26").Value = 44425
$ws.Range("B26").Value = 22
$ws.Range("C26").Value = 30.1

# Update the active selection to match where the author last clicked
[void]$ws.Range("C27").Select()
